$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 248, shifting existing rows 248:337 down to 250:339
$ws.Rows("248:249").Insert()

# Fill new row 248
$ws.Range("A248").Value = 9
$ws.Range("B248").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C248").Value = "Metropolitana"
$ws.Range("D248").Value = 44559
$ws.Range("E248").Value = 13
$ws.Range("F248").Value = 100112012
$ws.Range("G248").Value = "Espinaca"
$ws.Range("H248").Value = "Sin especificar"
$ws.Range("I248").Value = "Primera"
$ws.Range("J248").Value = 160
$ws.Range("K248").Value = 10000
$ws.Range("L248").Value = 12000
$ws.Range("M248").Value = 11000
$ws.Range("N248").Value = "`$/cuna 10 kilos"
$ws.Range("O248").Value = "Provincia de Chacabuco"
$ws.Range("P248").Value = 1100
$ws.Range("Q248").Value = 10
$ws.Range("R248").Value = "Hortaliza"

# Fill new row 249
$ws.Range("A249").Value = 9
$ws.Range("B249").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44559
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112012
$ws.Range("G249").Value = "Espinaca"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Segunda"
$ws.Range("J249").Value = 79
$ws.Range("K249").Value = 9000
$ws.Range("L249").Value = 9000
$ws.Range("M249").Value = 9000
$ws.Range("N249").Value = "`$/cuna 10 kilos"
$ws.Range("O249").Value = "Provincia de Chacabuco"
$ws.Range("P249").Value = 900
$ws.Range("Q249").Value = 10
$ws.Range("R249").Value = "Hortaliza"

# Ensure D column date formatting/number format matches the rest (numFmt YYYY-MM-DD HH:MM:SS)
$ws.Range("D248:D249").NumberFormat = $ws.Range("D250").NumberFormat
